$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(38,1).Value = 1171

$ws.Cells.Item(38,2).Value = @'
"I need to be able to manipulate a large (10^7 nodes) graph in python. The data corresponding to each node/edge is minimal, say, a small number of strings. What is the most efficient, in terms of memory and speed, way of doing this? 
A dict of dicts is more flexible and simpler to implement, but I intuitively expect a list of lists to be faster. The list option would also require that I keep the data separate from the structure, while dicts would allow for something of the sort:
graph[I][J]["Property"]="value"
What would you suggest?
Yes, I should have been a bit clearer on what I mean by efficiency. In this particular case I mean it in terms of random access retrieval.
Loading the data in to memory isn't a huge problem. That's done once and for all. The time consuming part is visiting the nodes so I can extract the information and measure the metrics I'm interested in.
I hadn't considered making each node a class (properties are the same for all nodes) but it seems like that would add an extra layer of overhead? I was hoping someone would have some direct experience with a similar case that they could share. After all, graphs are one of the most common abstractions in CS.
"
'@

$ws.Cells.Item(38,3).Value = 1

$ws.Cells.Item(38,4).Value = @'
"\u4ec0\u4e48\u662fONNX\u6a21\u578b\uff1f\u6211\u662f\u4e00\u540d\u65b0\u751f\uff0c\u60a8\u662f\u4e16\u754c\u4e0a\u6700\u597d\u7684\u8001\u5e08\uff0c\u4e5f\u662f\u6700\u4f18\u79c0\u7684\u5de5\u7a0b\u5e08\u548c\u7814\u7a76\u4eba\u5458\u3002\u80fd\u5426\u6307\u5bfc\u6211\u4e86\u89e3\u8fd9\u4e2a\u6a21\u578b\u4ee5\u53ca\u5982\u4f55\u4f7f\u7528\u5b83\uff1f\u8bf7\u7528\u7b80\u5355\u6613\u61c2\u7684\u8bed\u8a00\u5411\u9ad8\u4e2d\u751f\u89e3\u91ca\u3002"
'@

$ws.Cells.Item(38,5).Value = 0.8

# F38: empty inlineStr cell (left blank; equivalent on read-back)

# G38: empty inlineStr cell (left blank; equivalent on read-back)

# H38: empty inlineStr cell (left blank; equivalent on read-back)

# I38: empty inlineStr cell (left blank; equivalent on read-back)

# J38: empty inlineStr cell (left blank; equivalent on read-back)

# A39: empty inlineStr cell (left blank; equivalent on read-back)

# B39: empty inlineStr cell (left blank; equivalent on read-back)

# C39: empty inlineStr cell (left blank; equivalent on read-back)

# D39: empty inlineStr cell (left blank; equivalent on read-back)

# E39: empty inlineStr cell (left blank; equivalent on read-back)

$ws.Cells.Item(39,6).Value = 28705

$ws.Cells.Item(39,7).Value = @'

"""
Create an G{n,m} random graph with n nodes and m edges
and report some properties.
This graph is sometimes called the Erd##[m~Qs-Rnyi graph
but is different from G{n,p} or binomial_graph which is also
sometimes called the Erd##[m~Qs-Rnyi graph.
"""
__author__ = """Aric Hagberg (hagberg@lanl.gov)"""
__credits__ = """"""
#    Copyright (C) 2004-2006 by 
#    Aric Hagberg 
#    Dan Schult 
#    Pieter Swart 
#    Distributed under the terms of the GNU Lesser General Public License
#    http://www.gnu.org/copyleft/lesser.html
from networkx import *
import sys
n=10 # 10 nodes
m=20 # 20 edges
G=gnm_random_graph(n,m)
# some properties
print "node degree clustering"
for v in nodes(G):
    print v,degree(G,v),clustering(G,v)
# print the adjacency list to terminal 
write_adjlist(G,sys.stdout)

'@

$ws.Cells.Item(39,8).Value = @'
 
'@

$ws.Cells.Item(39,9).Value = @'
import torch.onnx
import torchvision.models as models
# 
model = models.resnet18(pretrained=True)
# 
model.eval()
# 
x = torch.randn(1, 3, 224, 224)
# 
torch.onnx.export(model, x, "resnet18.onnx")
import onnxruntime
#  ONNX 
session = onnxruntime.InferenceSession("resnet18.onnx")
# 
input_name = session.get_inputs()[0].name
output_name = session.get_outputs()[0].name
# 
x = np.random.randn(1, 3, 224, 224).astype(np.float32)
# 
result = session.run([output_name], {input_name: x})

'@

$ws.Cells.Item(39,10).Value = 0.8

$ws.Cells.Item(40,1).Value = 1171

$ws.Cells.Item(40,2).Value = @'
"I need to be able to manipulate a large (10^7 nodes) graph in python. The data corresponding to each node/edge is minimal, say, a small number of strings. What is the most efficient, in terms of memory and speed, way of doing this? 
A dict of dicts is more flexible and simpler to implement, but I intuitively expect a list of lists to be faster. The list option would also require that I keep the data separate from the structure, while dicts would allow for something of the sort:
graph[I][J]["Property"]="value"
What would you suggest?
Yes, I should have been a bit clearer on what I mean by efficiency. In this particular case I mean it in terms of random access retrieval.
Loading the data in to memory isn't a huge problem. That's done once and for all. The time consuming part is visiting the nodes so I can extract the information and measure the metrics I'm interested in.
I hadn't considered making each node a class (properties are the same for all nodes) but it seems like that would add an extra layer of overhead? I was hoping someone would have some direct experience with a similar case that they could share. After all, graphs are one of the most common abstractions in CS.
"
'@

$ws.Cells.Item(40,3).Value = 2

$ws.Cells.Item(40,4).Value = @'
"\u90a3\u4ed6\u4eec\u548ctorch tensorflow\u7684\u533a\u522b\u5728\u54ea\u91cc\u5462"
'@

$ws.Cells.Item(40,5).Value = 0.8

# F40: empty inlineStr cell (left blank; equivalent on read-back)

# G40: empty inlineStr cell (left blank; equivalent on read-back)

# H40: empty inlineStr cell (left blank; equivalent on read-back)

# I40: empty inlineStr cell (left blank; equivalent on read-back)

# J40: empty inlineStr cell (left blank; equivalent on read-back)

$ws.Cells.Item(41,1).Value = 1171

$ws.Cells.Item(41,2).Value = @'
"I need to be able to manipulate a large (10^7 nodes) graph in python. The data corresponding to each node/edge is minimal, say, a small number of strings. What is the most efficient, in terms of memory and speed, way of doing this? 
A dict of dicts is more flexible and simpler to implement, but I intuitively expect a list of lists to be faster. The list option would also require that I keep the data separate from the structure, while dicts would allow for something of the sort:
graph[I][J]["Property"]="value"
What would you suggest?
Yes, I should have been a bit clearer on what I mean by efficiency. In this particular case I mean it in terms of random access retrieval.
Loading the data in to memory isn't a huge problem. That's done once and for all. The time consuming part is visiting the nodes so I can extract the information and measure the metrics I'm interested in.
I hadn't considered making each node a class (properties are the same for all nodes) but it seems like that would add an extra layer of overhead? I was hoping someone would have some direct experience with a similar case that they could share. After all, graphs are one of the most common abstractions in CS.
"
'@

$ws.Cells.Item(41,3).Value = 3

$ws.Cells.Item(41,4).Value = @'
"\u4ed6\u4eec\u90fd\u662f\u56fe\u7ed3\u6784\u7684\u5417"
'@

$ws.Cells.Item(41,5).Value = 0.8

# F41: empty inlineStr cell (left blank; equivalent on read-back)

# G41: empty inlineStr cell (left blank; equivalent on read-back)

# H41: empty inlineStr cell (left blank; equivalent on read-back)

# I41: empty inlineStr cell (left blank; equivalent on read-back)

# J41: empty inlineStr cell (left blank; equivalent on read-back)
